# Update the year reference from 2024 to 2026 throughout the document.
# This mirrors the commit "fix: ubah tahun ke 2026" which changed the
# trailing "4" of "...202|4" to "6" in two places (turning "2024" -> "2026").

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "4",      # FindText
    $true,    # MatchCase
    $false,   # MatchWholeWord
    $false,   # MatchWildcards
    $false,   # MatchSoundsLike
    $false,   # MatchAllWordForms
    $true,    # Forward
    1,        # Wrap (wdFindContinue)
    $false,   # Format
    "6",      # ReplaceWith
    2         # Replace (wdReplaceAll)
)

$d.Save()
